$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.70682288476674
$ws.Range("C2").Value = 14.78532723153985
$ws.Range("D2").Value = 15.50347304566591
$ws.Range("E2").Value = 16.97685345109991
$ws.Range("G2").Value = 3.816270500864024
$ws.Range("J2").Value = 9.74281218610728
$ws.Range("K2").Value = 22.29088125580993
$ws.Range("N2").Value = 28.12527494638073
$ws.Range("B3").Value = 22.52684840088804
$ws.Range("C3").Value = 14.66039015836829
$ws.Range("D3").Value = 15.46400319587076
$ws.Range("E3").Value = 16.93919609860501
$ws.Range("G3").Value = 3.820759058267098
$ws.Range("J3").Value = 9.75865183894085
$ws.Range("K3").Value = 22.17230964206807
$ws.Range("N3").Value = 28.03750749653893
$ws.Range("B4").Value = 22.42209685708278
$ws.Range("C4").Value = 14.58754482269734
$ws.Range("D4").Value = 15.44346067881335
$ws.Range("E4").Value = 16.92013051928552
$ws.Range("G4").Value = 3.823653817088399
$ws.Range("J4").Value = 9.769951028963579
$ws.Range("K4").Value = 22.10511511426319
$ws.Range("N4").Value = 27.98449205057129
$ws.Range("B5").Value = 22.38089347738615
$ws.Range("C5").Value = 14.55885786120779
$ws.Range("D5").Value = 15.43602123218296
$ws.Range("E5").Value = 16.91338378643199
$ws.Range("G5").Value = 3.82486849804746
$ws.Range("J5").Value = 9.774950719251699
$ws.Range("K5").Value = 22.07916158309947
$ws.Range("N5").Value = 27.96311165059067
$ws.Range("B6").Value = 22.37414232239453
$ws.Range("C6").Value = 14.55415540966749
$ws.Range("D6").Value = 15.43484227618651
$ws.Range("E6").Value = 16.91232532124679
$ws.Range("G6").Value = 3.825072315796515
$ws.Range("J6").Value = 9.775804766042597
$ws.Range("K6").Value = 22.07493880689474
$ws.Range("N6").Value = 27.95957508224022
$ws.Range("B7").Value = 22.42153511968245
$ws.Range("C7").Value = 14.58715386641713
$ws.Range("D7").Value = 15.44335657123673
$ws.Range("E7").Value = 16.92003538697944
$ws.Range("G7").Value = 3.823670056612569
$ws.Range("J7").Value = 9.770016857165933
$ws.Range("K7").Value = 22.10475928812044
$ws.Range("N7").Value = 27.98420279563961
$ws.Range("B8").Value = 22.64359659862932
$ws.Range("C8").Value = 14.74146255626178
$ws.Range("D8").Value = 15.48909984377272
$ws.Range("E8").Value = 16.96302922852865
$ws.Range("G8").Value = 3.817789444604478
$ws.Range("J8").Value = 9.747946806259309
$ws.Range("K8").Value = 22.24884532512838
$ws.Range("N8").Value = 28.09483377468077
$ws.Range("B9").Value = 23.12295129832806
$ws.Range("C9").Value = 15.07350997442118
$ws.Range("D9").Value = 15.60790916980966
$ws.Range("E9").Value = 17.07936536099709
$ws.Range("G9").Value = 3.807351812542215
$ws.Range("J9").Value = 9.717175208212407
$ws.Range("K9").Value = 22.5749933413227
$ws.Range("N9").Value = 28.31853656661453
$ws.Range("B10").Value = 23.49936505607221
$ws.Range("C10").Value = 15.33364993664432
$ws.Range("D10").Value = 15.71263323333472
$ws.Range("E10").Value = 17.18409005223452
$ws.Range("G10").Value = 3.80034079753766
$ws.Range("J10").Value = 9.702221371148873
$ws.Range("K10").Value = 22.83983761912221
$ws.Range("N10").Value = 28.48676391459312
$ws.Range("B11").Value = 23.67528699854266
$ws.Range("C11").Value = 15.45510325295286
$ws.Range("D11").Value = 15.76397594706408
$ws.Range("E11").Value = 17.23583509790392
$ws.Range("G11").Value = 3.797292027253155
$ws.Range("J11").Value = 9.697085820588114
$ws.Range("K11").Value = 22.9654718044922
$ws.Range("N11").Value = 28.56407919756564
$ws.Range("B12").Value = 23.74252738120184
$ws.Range("C12").Value = 15.50150672747201
$ws.Range("D12").Value = 15.78394199056003
$ws.Range("E12").Value = 17.25601182678704
$ws.Range("G12").Value = 3.796157592954246
$ws.Range("J12").Value = 9.695381165587994
$ws.Range("K12").Value = 23.01375740043413
$ws.Range("N12").Value = 28.59346434702127
$ws.Range("B13").Value = 23.72801910823867
$ws.Range("C13").Value = 15.49149517838976
$ws.Range("D13").Value = 15.77961880507787
$ws.Range("E13").Value = 17.25164066451581
$ws.Range("G13").Value = 3.796401023430088
$ws.Range("J13").Value = 9.695737611525033
$ws.Range("K13").Value = 23.00332712955278
$ws.Range("N13").Value = 28.58713104079091
$ws.Range("B14").Value = 23.68080672843816
$ws.Range("C14").Value = 15.45891285195097
$ws.Range("D14").Value = 15.76560812392954
$ws.Range("E14").Value = 17.23748343617926
$ws.Range("G14").Value = 3.797198295285866
$ws.Range("J14").Value = 9.696940764493721
$ws.Range("K14").Value = 22.96943022744369
$ws.Range("N14").Value = 28.56649461018467
$ws.Range("B15").Value = 23.65196730407672
$ws.Range("C15").Value = 15.43900777412083
$ws.Range("D15").Value = 15.75709409799788
$ws.Range("E15").Value = 17.22888725667905
$ws.Range("G15").Value = 3.797689256548967
$ws.Range("J15").Value = 9.697709003635859
$ws.Range("K15").Value = 22.94875901094924
$ws.Range("N15").Value = 28.55386802957985
$ws.Range("B16").Value = 23.48795810595022
$ws.Range("C16").Value = 15.32577227553065
$ws.Range("D16").Value = 15.7093516474124
$ws.Range("E16").Value = 17.18079030865715
$ws.Range("G16").Value = 3.800542858123086
$ws.Range("J16").Value = 9.702590567464489
$ws.Range("K16").Value = 22.83172827752752
$ws.Range("N16").Value = 28.48172654467838
$ws.Range("B17").Value = 23.38850837477922
$ws.Range("C17").Value = 15.25707817136103
$ws.Range("D17").Value = 15.68100565683901
$ws.Range("E17").Value = 17.15233020460516
$ws.Range("G17").Value = 3.802329351293044
$ws.Range("J17").Value = 9.70601248050165
$ws.Range("K17").Value = 22.7612327322256
$ws.Range("N17").Value = 28.4376678318511
$ws.Range("B18").Value = 23.3317520922778
$ws.Range("C18").Value = 15.21786249070168
$ws.Range("D18").Value = 15.66505094093848
$ws.Range("E18").Value = 17.13634766861273
$ws.Range("G18").Value = 3.80337013526571
$ws.Range("J18").Value = 9.708137565196351
$ws.Range("K18").Value = 22.72117272256357
$ws.Range("N18").Value = 28.41240093738706
$ws.Range("B19").Value = 23.31261330520515
$ws.Range("C19").Value = 15.20463655000051
$ws.Range("D19").Value = 15.65970917048168
$ws.Range("E19").Value = 17.131002956842
$ws.Range("G19").Value = 3.803724805318738
$ws.Range("J19").Value = 9.708884016919017
$ws.Range("K19").Value = 22.70769367277171
$ws.Range("N19").Value = 28.40385897037399
$ws.Range("B20").Value = 23.39904932850778
$ws.Range("C20").Value = 15.26436046759992
$ws.Range("D20").Value = 15.68398706183376
$ws.Range("E20").Value = 17.1553198395463
$ws.Range("G20").Value = 3.802137806776133
$ws.Range("J20").Value = 9.705631971704619
$ws.Range("K20").Value = 22.76868690135762
$ws.Range("N20").Value = 28.44235028638279
$ws.Range("B21").Value = 23.69465767038204
$ws.Range("C21").Value = 15.46847218747421
$ws.Range("D21").Value = 15.76970926291414
$ws.Range("E21").Value = 17.24162603269555
$ws.Range("G21").Value = 3.79696357371482
$ws.Range("J21").Value = 9.696580851445464
$ws.Range("K21").Value = 22.97936753343641
$ws.Range("N21").Value = 28.57255315869075
$ws.Range("B22").Value = 23.89145891201392
$ws.Range("C22").Value = 15.60425377888233
$ws.Range("D22").Value = 15.82878135833665
$ws.Range("E22").Value = 17.30141999798415
$ws.Range("G22").Value = 3.793698822277772
$ws.Range("J22").Value = 9.692064836119643
$ws.Range("K22").Value = 23.12118659790806
$ws.Range("N22").Value = 28.65827221015078
$ws.Range("B23").Value = 23.78610996142186
$ws.Range("C23").Value = 15.53157864558112
$ws.Range("D23").Value = 15.79697774100915
$ws.Range("E23").Value = 17.26919982085729
$ws.Range("G23").Value = 3.795430632513912
$ws.Range("J23").Value = 9.694346970414342
$ws.Range("K23").Value = 23.0451279415916
$ws.Range("N23").Value = 28.61246720422449
$ws.Range("B24").Value = 23.39428245527884
$ws.Range("C24").Value = 15.26106727446794
$ws.Range("D24").Value = 15.6826381034324
$ws.Range("E24").Value = 17.15396704273304
$ws.Range("G24").Value = 3.80222436135266
$ws.Range("J24").Value = 9.705803508305108
$ws.Range("K24").Value = 22.7653154089205
$ws.Range("N24").Value = 28.44023315187917
$ws.Range("B25").Value = 22.98882623147512
$ws.Range("C25").Value = 14.98070616515306
$ws.Range("D25").Value = 15.5726802760198
$ws.Range("E25").Value = 17.0444889528632
$ws.Range("G25").Value = 3.810059312501354
$ws.Range("J25").Value = 9.72415751159008
$ws.Range("K25").Value = 22.48222645846734
$ws.Range("N25").Value = 28.25732014719682
